# Insert a new weekly price record as the new first row of the Papaya /
# Vega Modelo de Temuco data block (row 41), pushing the existing rows
# 41-62 down to 42-63. Mirrors Excel's "Insert Sheet Rows" behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a blank row at row 41.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A41").Value2 = 10
$ws.Range("B41").Value2 = "Vega Modelo de Temuco"
$ws.Range("C41").Value2 = "La Araucanía"
$ws.Range("D41").Value2 = 44489
$ws.Range("E41").Value2 = 9
$ws.Range("F41").Value2 = "Fruta"
$ws.Range("G41").Value2 = 100108
$ws.Range("H41").Value2 = "Tropicales y subtropicales"
$ws.Range("I41").Value2 = 100108004
$ws.Range("J41").Value2 = "Papaya"
$ws.Range("K41").Value2 = "Cultivar IV Región"
$ws.Range("L41").Value2 = "Primera"
$ws.Range("M41").Value2 = 100
$ws.Range("N41").Value2 = 20000
$ws.Range("O41").Value2 = 20000
$ws.Range("P41").Value2 = 20000
$ws.Range("Q41").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R41").Value2 = "Provincia del Elquí"
$ws.Range("S41").Value2 = 2000
$ws.Range("T41").Value2 = 10

# Make sure the date cell keeps the workbook's date number format (style
# index 2 in styles.xml), matching the other date cells in column D.
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
